$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be misread as a number by Excel (e.g. "596.53")
# need the cell pre-formatted as Text so the literal string is preserved,
# exactly like the original file (all cells are plain text/inlineStr).
$textCells = @(
    'D5', 'D6', 'D8', 'D12', 'D13', 'D16', 'D19', 'D21', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '74.734.68'
$ws.Range('E2').Value = '  -0.18%  '

$ws.Range('D3').Value = '2.817.58'
$ws.Range('E3').Value = '  +8.06%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '596.53'
$ws.Range('E5').Value = '  +2.38%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = '187.37'
$ws.Range('E6').Value = '  +0.29%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '0.557'
$ws.Range('E8').Value = '  +3.42%  '

$ws.Range('E9').Value = '  -8.19%  '

$ws.Range('D10').Value = '2.820.03'
$ws.Range('E10').Value = '  +8.03%  '

$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').Value = '0.371'
$ws.Range('E12').Value = '  +1.84%  '

$ws.Range('D13').Value = '4.84'
$ws.Range('E13').Value = '  +0.48%  '

$ws.Range('D14').Value = '3.342.42'
$ws.Range('E14').Value = '  +8.17%  '

$ws.Range('D15').Value = '74.831.40'
$ws.Range('E15').Value = '  -0.31%  '

$ws.Range('D16').Value = '27.06'
$ws.Range('E16').Value = '  +2.68%  '

$ws.Range('E17').Value = '  -3.63%  '

$ws.Range('D18').Value = '2.830.61'
$ws.Range('E18').Value = '  +7.98%  '

$ws.Range('D19').Value = '8.93'
$ws.Range('E19').Value = '  -3.42%  '

$ws.Range('E20').Value = '  +4.55%  '

$ws.Range('D21').Value = '374.63'
$ws.Range('E21').Value = '  -1.15%  '

$ws.Range('E22').Value = '  -2.89%  '

$ws.Range('D23').Value = '4.11'
$ws.Range('E23').Value = '  +0.52%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '70.73'
$ws.Range('E25').Value = '  +0.72%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.971.39'
$ws.Range('E26').Value = '  +8.39%  '

$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = '4.18'
$ws.Range('E27').Value = '  -0.87%  '

$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').Value = '9.58'
$ws.Range('E28').Value = '  +2.29%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0000103'
$ws.Range('E29').Value = '  +8.24%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '519.82'
$ws.Range('E31').Value = '  +1.59%  '

$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  -0.63%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '7.89'
$ws.Range('E33').Value = '  -0.89%  '

$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').Value = '1.80'
$ws.Range('E34').Value = '  +2.97%  '

$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').Value = '20.14'
$ws.Range('E36').Value = '  +4.72%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '162.69'
$ws.Range('E37').Value = '  +2.26%  '

$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -1.20%  '

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '19.28'
$ws.Range('E39').Value = '  -0.59%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '182.42'
$ws.Range('E40').Value = '  +16.45%  '

$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.07%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '5.05'
$ws.Range('E42').Value = '  +1.95%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '0.339'
$ws.Range('E43').Value = '  +3.78%  '

$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '1.68'
$ws.Range('E44').Value = '  -1.38%  '

$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').Value = '1.23'
$ws.Range('E45').Value = '  +4.07%  '

$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '39.74'
$ws.Range('E46').Value = '  +2.29%  '

$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '2.34'
$ws.Range('E47').Value = '  -5.29%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0857'
$ws.Range('E48').Value = '  +3.02%  '

$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '0.568'
$ws.Range('E49').Value = '  +7.73%  '

$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '3.74'
$ws.Range('E50').Value = '  +2.92%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.615'
$ws.Range('E51').Value = '  +4.83%  '

# Restore default (General) style on the forced-text cells so no residual
# cell-level style/number-format difference remains versus the original.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}